$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.329955947548125
$ws.Range("B2").Value = 0.3835379738755286
$ws.Range("C2").Value = 0.2821739802348961
$ws.Range("D2").Value = 0.3549827534312814
$ws.Range("E2").Value = 0.3038799051185036
$ws.Range("K2").Value = 2.309691632836875
$ws.Range("L2").Value = 2.6847658171287
$ws.Range("M2").Value = 1.975217861644273
$ws.Range("N2").Value = 2.48487927401897
$ws.Range("O2").Value = 2.127159335829525
$ws.Range("P2").Value = 176.00503
$ws.Range("Q2").Value = 234.6407409430736
$ws.Range("R2").Value = 134.4303197968127
$ws.Range("S2").Value = 198.927849887097
$ws.Range("T2").Value = 152.3979670103585
$ws.Range("U2").Value = 0.04045397268321832
$ws.Range("V2").Value = 0.06075139190588973
$ws.Range("W2").Value = 0.021830316360436
$ws.Range("X2").Value = 0.05074990909545279
$ws.Range("Y2").Value = 0.03001347679060722
$ws.Range("Z2").Value = 0.6761959217999651
$ws.Range("AA2").Value = 0.7867371493239307
$ws.Range("AB2").Value = 0.5428695047209194
$ws.Range("AC2").Value = 0.7361708161371121
$ws.Range("AD2").Value = 0.6077257788357813
$ws.Range("F3").Value = 6.995053518910998
$ws.Range("G3").Value = 8.839945533553912
$ws.Range("H3").Value = 5.458409058833954
$ws.Range("I3").Value = 7.827740952391174
$ws.Range("J3").Value = 6.128389689880177
$ws.Range("K3").Value = 2.308367661240629
$ws.Range("L3").Value = 2.917182026072791
$ws.Range("M3").Value = 1.801274989415205
$ws.Range("N3").Value = 2.583154514289088
$ws.Range("O3").Value = 2.022368597660459
$ws.Range("P3").Value = 171.71772
$ws.Range("Q3").Value = 198.287416170475
$ws.Range("R3").Value = 150.1504258470239
$ws.Range("S3").Value = 182.3819444923535
$ws.Range("T3").Value = 159.821941354878
$ws.Range("U3").Value = 0.04073408270762569
$ws.Range("V3").Value = 0.0463885921698482
$ws.Range("W3").Value = 0.03518715254467193
$ws.Range("X3").Value = 0.04365062156813727
$ws.Range("Y3").Value = 0.03769072520242973
$ws.Range("Z3").Value = 0.6835048984848644
$ws.Range("AA3").Value = 0.6835064434634854
$ws.Range("AB3").Value = 0.6835002713780689
$ws.Range("AC3").Value = 0.6835052668306218
$ws.Range("AD3").Value = 0.6835037161834524
$ws.Range("P4").Value = 186.11938
$ws.Range("Q4").Value = 331.1211714758743
$ws.Range("R4").Value = 124.8773891949485
$ws.Range("S4").Value = 225.9438365158736
$ws.Range("T4").Value = 143.8209948344143
$ws.Range("U4").Value = 0.04217750107223983
$ws.Range("V4").Value = 0.09856727039136116
$ws.Range("W4").Value = 0.002995157064678238
$ws.Range("X4").Value = 0.06900750304990227
$ws.Range("Y4").Value = 0.0179721845264825
$ws.Range("Z4").Value = 0.6480967536251064
$ws.Range("AA4").Value = 0.8788590562026701
$ws.Range("AB4").Value = 0.2942739160203218
$ws.Range("AC4").Value = 0.7754857083400452
$ws.Range("AD4").Value = 0.4845112728795999
$ws.Range("A5").Value = 0.3297527023369954
$ws.Range("B5").Value = 0.3835211989853209
$ws.Range("C5").Value = 0.2820646983561017
$ws.Range("D5").Value = 0.3548358448971214
$ws.Range("E5").Value = 0.3036844575333699
$ws.Range("F5").Value = 6.99959957331527
$ws.Range("G5").Value = 8.822211548714666
$ws.Range("H5").Value = 5.468825961423108
$ws.Range("I5").Value = 7.833020355743384
$ws.Range("J5").Value = 6.128967116855383
$ws.Range("K5").Value = 2.308174082155408
$ws.Range("L5").Value = 3.036285692548471
$ws.Range("M5").Value = 1.71589277237482
$ws.Range("N5").Value = 2.633113997229116
$ws.Range("O5").Value = 1.967495510641607
$ws.Range("P5").Value = 187.77498
$ws.Range("Q5").Value = 363.061766671988
$ws.Range("R5").Value = 108.7124494225298
$ws.Range("S5").Value = 237.8216208933891
$ws.Range("T5").Value = 137.053414814629
$ws.Range("U5").Value = 0.04215473926513635
$ws.Range("V5").Value = 0.103722803249861
$ws.Range("W5").Value = 0.0009062563108073947
$ws.Range("X5").Value = 0.07112792526370111
$ws.Range("Y5").Value = 0.01641321719512654
$ws.Range("Z5").Value = 0.6380592361797547
$ws.Range("AA5").Value = 0.8940501986107428
$ws.Range("AB5").Value = 0.2299770683133065
$ws.Range("AC5").Value = 0.7798247068670745
$ws.Range("AD5").Value = 0.4523311835122468
